$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 374, shifting existing rows 374:394 down to 375:395
$ws.Rows("374:374").Insert()

# Populate the newly inserted row 374 with the new data record
$ws.Range("A374").Value = 6
$ws.Range("B374").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C374").Value = "Metropolitana"
$ws.Range("D374").Value = 45267
$ws.Range("E374").Value = 13
$ws.Range("F374").Value = "Fruta"
$ws.Range("G374").Value = 100101
$ws.Range("H374").Value = "Berries"
$ws.Range("I374").Value = 100101004
$ws.Range("J374").Value = "Frambuesa"
$ws.Range("K374").Value = "Sin especificar"
$ws.Range("L374").Value = "Primera"
$ws.Range("M374").Value = 1000
$ws.Range("N374").Value = 13000
$ws.Range("O374").Value = 13000
$ws.Range("P374").Value = 13000
$ws.Range("Q374").Value = "`$/bandeja 2 kilos"
$ws.Range("R374").Value = "Provincia de Curicó"
$ws.Range("S374").Value = 6500
$ws.Range("T374").Value = 2
